# Editar Produto e Finalizar Venda pronto
#
# This script reproduces, via Excel COM automation, the edits that were made to
# the "Compras.xlsx" workbook:
#   1) The purchase "Método" (method) catalog is updated: the existing method
#      "Pacote" is renamed to "Unidade" and a brand new method "Combo" is added.
#   2) The single product record is edited from "Ração / Pedigree / Pacote" to
#      "Churrasqueira / Grande / Combo", with updated purchase/sale prices.
#   3) A sale is staged (P_Vendas helper table) and then finalized (Vendas
#      table), and a purchase is staged (P_Compras helper table) and finalized
#      (Compras table), updating the Estoque (stock) table accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Métodos sheet - rename "Pacote" -> "Unidade" and add new method "Combo"
# ---------------------------------------------------------------------------
$wsMetodos = $wb.Worksheets.Item("Métodos")
$wsMetodos.Range("A2").Value = "Unidade"
$wsMetodos.Range("A3").Value = "Combo"

# ---------------------------------------------------------------------------
# 2) Produtos sheet - edit the product record in place
# ---------------------------------------------------------------------------
$wsProdutos = $wb.Worksheets.Item("Produtos")
$wsProdutos.Range("A2").Value = "Churrasqueira"
$wsProdutos.Range("B2").Value = "Grande"
$wsProdutos.Range("C2").Value = "Combo"
$wsProdutos.Range("D2").Value = "Combo"
# Valor_Compra / Valor_Venda are stored as text ("70"/"100") in the source
# data, so force a text number format before assigning to avoid Excel's
# automatic numeric coercion of digit-only strings.
$wsProdutos.Range("E2").NumberFormat = "@"
$wsProdutos.Range("E2").Value = "70"
$wsProdutos.Range("F2").NumberFormat = "@"
$wsProdutos.Range("F2").Value = "100"
$wsProdutos.Range("G2").Value = "Não"

# ---------------------------------------------------------------------------
# 3) P_Vendas sheet - stage the sale cart lines
# ---------------------------------------------------------------------------
$wsPVendas = $wb.Worksheets.Item("P_Vendas")
$wsPVendas.Range("A2").Value = 10000
$wsPVendas.Range("B2").Value = "Churrasqueira"
$wsPVendas.Range("C2").Value = "Grande"
$wsPVendas.Range("D2").Value = "Combo"
$wsPVendas.Range("E2").Value = 10
$wsPVendas.Range("F2").Value = 100
$wsPVendas.Range("G2").Value = 1000

$wsPVendas.Range("A3").Value = 10000
$wsPVendas.Range("B3").Value = "Churrasqueira"
$wsPVendas.Range("C3").Value = "Grande"
$wsPVendas.Range("D3").Value = "Combo"
$wsPVendas.Range("E3").Value = 20
$wsPVendas.Range("F3").Value = 100
$wsPVendas.Range("G3").Value = 2000

$wsPVendas.Range("A4").Value = 10000
$wsPVendas.Range("B4").Value = "Churrasqueira"
$wsPVendas.Range("C4").Value = "Grande"
$wsPVendas.Range("D4").Value = "Combo"
$wsPVendas.Range("E4").Value = 20
$wsPVendas.Range("F4").Value = 100
$wsPVendas.Range("G4").Value = 2000

# ---------------------------------------------------------------------------
# 4) Vendas sheet - finalize the sale (QItens, Frete?, Desconto, Pagamento, ...)
# ---------------------------------------------------------------------------
$wsVendas = $wb.Worksheets.Item("Vendas")
$wsVendas.Range("A2").Value = 10000
$wsVendas.Range("B2").Value = 50
$wsVendas.Range("C2").Value = "Sim"
$wsVendas.Range("D2").Value = 953
$wsVendas.Range("E2").Value = "Pix"
$wsVendas.Range("G2").Value = 4047
$wsVendas.Range("H2").Value = "15/01/2023"

# ---------------------------------------------------------------------------
# 5) P_Compras sheet - stage the purchase line
# ---------------------------------------------------------------------------
$wsPCompras = $wb.Worksheets.Item("P_Compras")
$wsPCompras.Range("A2").Value = 1
$wsPCompras.Range("B2").Value = "Churrasqueira"
$wsPCompras.Range("C2").Value = "Grande"
$wsPCompras.Range("D2").Value = "Combo"
$wsPCompras.Range("E2").Value = 100
$wsPCompras.Range("F2").Value = 70
$wsPCompras.Range("G2").Value = 7000

# ---------------------------------------------------------------------------
# 6) Compras sheet - finalize the purchase
# ---------------------------------------------------------------------------
$wsCompras = $wb.Worksheets.Item("Compras")
$resumo = "ID                         1" + "`n" +
          "Produto        Churrasqueira" + "`n" +
          "Marca               Pedigree" + "`n" +
          "Método               Unidade" + "`n" +
          "Quantidade               100" + "`n" +
          "Valor_Un                70.0" + "`n" +
          "Valor_Total           7000.0" + "`n" +
          "dtype: object"
$wsCompras.Range("B2").Value = $resumo
$wsCompras.Range("C2").Value = "15/01/2023"
$wsCompras.Range("D2").Value = 7000

# ---------------------------------------------------------------------------
# 7) Estoque sheet - updated stock quantity for the (edited) product
# ---------------------------------------------------------------------------
$wsEstoque = $wb.Worksheets.Item("Estoque")
$wsEstoque.Range("A2").Value = "Churrasqueira"
$wsEstoque.Range("B2").Value = "Grande"
$wsEstoque.Range("C2").Value = "Combo"
$wsEstoque.Range("D2").Value = 50
